{"js": "// Apply two anonymization edits to the loan contract:\n// 1) \"V\u011b\u0159itel\u00e9 poskytuj\u00ed dlu\u017en\u00edk\u016fm spole\u010dnou p\u016fj\u010dku ve v\u00fd\u0161i<br/>sto dvacet\n//    tis\u00edc korun \u010desk\u00fdch (120 000 K\u010d).\" -> replace the literal amount\n//    \"120 000 K\u010d\" with the placeholder \"[[AMOUNT_1]]\". The run that used to\n//    carry the bold amount is merged with its neighbours into a single,\n//    non-bold run (the whole paragraph is reset to plain text, preserving\n//    the manual line break between the two sentences).\n// 2) \"...smluvn\u00ed pokutu ve v\u00fd\u0161i 1 000 K\u010d za ka\u017ed\u00fd zapo\u010dat\u00fd m\u011bs\u00edc\n//    prodlen\u00ed.\" -> replace \"1 000 K\u010d\" with \"[[AMOUNT_2]]\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// --- Edit 1: collapse the \"p\u016fj\u010dku ve v\u00fd\u0161i\" paragraph into one plain run ---\nlet loanParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"V\u011b\u0159itel\u00e9 poskytuj\u00ed dlu\u017en\u00edk\u016fm spole\u010dnou p\u016fj\u010dku ve v\u00fd\u0161i\") !== -1) {\n    loanParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!loanParagraph) {\n  throw new Error(\"Could not find the 'V\u011b\u0159itel\u00e9 poskytuj\u00ed...' paragraph\");\n}\n\n// \"\\v\" (vertical tab, 0x0B) is how Office.js represents a manual line break\n// (<w:br/>) inside a text string passed to insertText.\nloanParagraph.insertText(\n  \"V\u011b\u0159itel\u00e9 poskytuj\u00ed dlu\u017en\u00edk\u016fm spole\u010dnou p\u016fj\u010dku ve v\u00fd\u0161i\\vsto dvacet tis\u00edc korun \u010desk\u00fdch ([[AMOUNT_1]]).\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- Edit 2: swap the contractual-penalty amount for its placeholder ---\nconst penaltyAmount = body.search(\"1 000 K\u010d\", { matchCase: true, matchWholeWord: false });\npenaltyAmount.load(\"text\");\nawait context.sync();\n\nif (penaltyAmount.items.length === 0) {\n  throw new Error(\"Could not find '1 000 K\u010d' to replace\");\n}\n\npenaltyAmount.items[0].insertText(\"[[AMOUNT_2]]\", \"Replace\");\nawait context.sync();\n", "ps1": "# Apply two anonymization edits to the loan contract:\n# 1) \"Veritele poskytuji dluznikum spolecnou pujcku ve vysi<br/>sto dvacet\n#    tisic korun ceskych (120 000 Kc).\" -> replace the literal amount\n#    \"120 000 Kc\" with the placeholder \"[[AMOUNT_1]]\". The run that used to\n#    carry the bold amount is merged with its neighbours into a single,\n#    non-bold run (the whole paragraph is rewritten as plain text, keeping\n#    the manual line break between the two sentences).\n# 2) \"...smluvni pokutu ve vysi 1 000 Kc za kazdy zapocaty mesic\n#    prodleni.\" -> replace \"1 000 Kc\" with \"[[AMOUNT_2]]\".\n\n$d = $word.ActiveDocument\n\n# Vertical tab (chr 11) is how Word represents a manual line break\n# (<w:br/>) inside Range.Text.\n$vt = [char]11\n\n# --- Edit 1: collapse the \"pujcku ve vysi\" paragraph into one plain run ---\n$loanParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*V\u011b\u0159itel\u00e9 poskytuj\u00ed dlu\u017en\u00edk\u016fm spole\u010dnou p\u016fj\u010dku ve v\u00fd\u0161i*\") {\n        $loanParagraph = $p\n        break\n    }\n}\n\nif ($loanParagraph -eq $null) {\n    throw \"Could not find the 'V\u011b\u0159itel\u00e9 poskytuj\u00ed...' paragraph\"\n}\n\n$r = $loanParagraph.Range\n# Exclude the trailing paragraph mark so we don't merge with the next paragraph.\n$sub = $d.Range($r.Start, $r.End - 1)\n$sub.Text = \"V\u011b\u0159itel\u00e9 poskytuj\u00ed dlu\u017en\u00edk\u016fm spole\u010dnou p\u016fj\u010dku ve v\u00fd\u0161i\" + $vt + \"sto dvacet tis\u00edc korun \u010desk\u00fdch ([[AMOUNT_1]]).\"\n\n# --- Edit 2: swap the contractual-penalty amount for its placeholder ---\n$find = $d.Content.Find\n$find.Text = \"1 000 K\u010d\"\n$find.Replacement.Text = \"[[AMOUNT_2]]\"\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    throw \"Could not find '1 000 K\u010d' to replace\"\n}\n"}
